$d = $word.ActiveDocument

$replacements = @(
    @("39×48=", "36×53="),
    @("99×66=", "73×51="),
    @("44×56=", "73×93="),
    @("87×82=", "27×81="),
    @("28×15=", "41×21="),
    @("73×96=", "87×79="),
    @("41×37=", "26×80="),
    @("94×48=", "23×75="),
    @("74×49=", "97×94="),
    @("29×86=", "27×47="),
    @("17×58=", "97×11="),
    @("92×70=", "31×88="),
    @("76×98=", "66×86="),
    @("85×90=", "21×25="),
    @("32×92=", "46×36="),
    @("12×99=", "92×12="),
    @("13×70=", "43×12="),
    @("94×96=", "75×44="),
    @("26×45=", "38×15="),
    @("72×58=", "39×43="),
    @("97×69=", "91×56="),
    @("55×30=", "69×26="),
    @("30×25=", "43×66="),
    @("66×44=", "45×56="),
    @("15×59=", "93×93=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
